# The "Output:" code-block paragraph currently reads (across three runs,
# separated by manual line breaks):
#   "Student@Student:~/Desktop/60004170098$" <br> "./sample" <br> "Hello World" <br> "Lucifer"
#
# It needs to become (merged into a single run, removing the line break
# between the prompt and the command, and updating the roll number /
# adding the ".py" extension):
#   "Student@Student:~/Desktop/60004170081$ ./sample.py" <br> "Hello World" <br> "Lucifer"
#
# Find the paragraph that contains the shell-prompt text so we don't rely
# on a hard-coded paragraph index.
$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Student@Student:~/Desktop/60004170098`$*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $paraStart = $target.Range.Start

    # Locate the old prompt text and the following "./sample" text (which
    # are split across two runs joined by a manual line break) so we can
    # replace that whole span - including the line break between them -
    # with a single merged run.
    $promptRange = $d.Content
    $promptRange.Start = $paraStart
    $found = $promptRange.Find.Execute("Student@Student:~/Desktop/60004170098`$", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($found -and $promptRange.Start -eq $paraStart) {
        $afterPrompt = $promptRange.End

        $sampleRange = $d.Range($afterPrompt, $target.Range.End)
        $sampleRange.Find.Execute("./sample", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $afterSample = $sampleRange.End

        # Replace from the start of the prompt through the end of
        # "./sample" (spanning the intervening run break) with the new,
        # single merged run of text.
        $replaceRange = $d.Range($paraStart, $afterSample)
        $replaceRange.Text = "Student@Student:~/Desktop/60004170081`$ ./sample.py"
    }
}
